# Auto-generated edit script: updates market price / profit columns (H-N)
# in the per-job Leve tables, mirroring a scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 13
$ws.Range("H13").Value = 5966.6665
$ws.Range("I13").Value = 3950
$ws.Range("J13").Value = 10000
$ws.Range("K13").Value = 3950
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = -3781
$ws.Range("N13").Value = -10338
# Row 98
$ws.Range("H98").Value = 45934.5
$ws.Range("I98").Value = 46924
$ws.Range("K98").Value = 46924
$ws.Range("M98").Value = -45426
# Row 122
$ws.Range("H122").Value = 45934.5
$ws.Range("I122").Value = 46924
$ws.Range("K122").Value = 140772
$ws.Range("M122").Value = -138322
# Row 138
$ws.Range("H138").Value = 2988.4915
$ws.Range("I138").Value = 1197.5312
$ws.Range("K138").Value = 3592.5936
$ws.Range("M138").Value = 1547.4064

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2853
$ws.Range("I32").Value = 2853
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2853
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -2566
# Row 135
$ws.Range("H135").Value = 105198.14
$ws.Range("J135").Value = 105198.14
$ws.Range("L135").Value = 105198.14
$ws.Range("N135").Value = -115338.14

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 32
$ws.Range("H32").Value = 52977.75
$ws.Range("J32").Value = 52977.75
$ws.Range("L32").Value = 52977.75
$ws.Range("N32").Value = -53745.75
# Row 39
$ws.Range("H39").Value = 30041
$ws.Range("J39").Value = 30041
$ws.Range("L39").Value = 30041
$ws.Range("N39").Value = -30819
# Row 118
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").ClearContents()
$ws.Range("N118").Value = 0
# Row 138
$ws.Range("H138").Value = 156486.88
$ws.Range("J138").Value = 156486.88
$ws.Range("L138").Value = 156486.88
$ws.Range("N138").Value = -166766.88

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 21
$ws.Range("H21").Value = 6300.5557
$ws.Range("I21").Value = 3398
$ws.Range("J21").Value = 7129.857
$ws.Range("K21").Value = 3398
$ws.Range("L21").Value = 7129.857
$ws.Range("M21").Value = -3163
$ws.Range("N21").Value = -7599.857
# Row 32
$ws.Range("H32").Value = 4287013.5
$ws.Range("I32").Value = 4287013.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4287013.5
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -4286697.5
# Row 38
$ws.Range("H38").Value = 1332.6666
$ws.Range("I38").Value = 1332.6666
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 1332.6666
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -955.6666
# Row 39
$ws.Range("H39").Value = 4910.2
$ws.Range("I39").Value = 4887.75
$ws.Range("J39").Value = 5000
$ws.Range("K39").Value = 4887.75
$ws.Range("L39").Value = 5000
$ws.Range("M39").Value = -4496.75
$ws.Range("N39").Value = -5782
# Row 46
$ws.Range("H46").Value = 1332.6666
$ws.Range("I46").Value = 1332.6666
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1332.6666
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -1121.6666
# Row 49
$ws.Range("H49").Value = 4910.2
$ws.Range("I49").Value = 4887.75
$ws.Range("J49").Value = 5000
$ws.Range("K49").Value = 4887.75
$ws.Range("L49").Value = 5000
$ws.Range("M49").Value = -4705.75
$ws.Range("N49").Value = -5364
# Row 132
$ws.Range("H132").Value = 22294.55
$ws.Range("I132").Value = 1932.2222
$ws.Range("K132").Value = 5796.6666
$ws.Range("M132").Value = -3266.6666
# Row 134
$ws.Range("H134").Value = 3384.0852
$ws.Range("I134").Value = 986.8333
$ws.Range("K134").Value = 2960.4999
$ws.Range("M134").Value = -425.4998999999998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 1165.7778
$ws.Range("I17").Value = 500.2857
$ws.Range("J17").Value = 3495
$ws.Range("K17").Value = 1500.8571
$ws.Range("L17").Value = 10485
$ws.Range("M17").Value = -1331.8571
$ws.Range("N17").Value = -10823
# Row 56
$ws.Range("H56").Value = 5378.727
$ws.Range("I56").Value = 5378.727
$ws.Range("K56").Value = 5378.727
$ws.Range("M56").Value = -4848.727
# Row 98
$ws.Range("H98").Value = 783
$ws.Range("I98").Value = 974.7
$ws.Range("K98").Value = 2924.1
$ws.Range("M98").Value = -1426.1
# Row 131
$ws.Range("H131").Value = 5182.1577
$ws.Range("I131").Value = 5279.125
$ws.Range("K131").Value = 15837.375
$ws.Range("M131").Value = -10797.375

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 32999.8
$ws.Range("J57").Value = 44999
$ws.Range("L57").Value = 44999
$ws.Range("N57").Value = -46639
# Row 70
$ws.Range("H70").Value = 6444.7085
$ws.Range("I70").Value = 5526.625
$ws.Range("K70").Value = 5526.625
$ws.Range("M70").Value = -5256.625
# Row 73
$ws.Range("H73").Value = 6444.7085
$ws.Range("I73").Value = 5526.625
$ws.Range("K73").Value = 5526.625
$ws.Range("M73").Value = -4590.625
# Row 80
$ws.Range("H80").Value = 21166.834
$ws.Range("I80").Value = 25764.572
$ws.Range("J80").Value = 14730
$ws.Range("K80").Value = 25764.572
$ws.Range("L80").Value = 14730
$ws.Range("M80").Value = -24766.572
$ws.Range("N80").Value = -16726
# Row 83
$ws.Range("H83").Value = 21166.834
$ws.Range("I83").Value = 25764.572
$ws.Range("J83").Value = 14730
$ws.Range("K83").Value = 128822.86
$ws.Range("L83").Value = 73650
$ws.Range("M83").Value = -123830.86
$ws.Range("N83").Value = -83634
# Row 102
$ws.Range("H102").Value = 4089.4897
$ws.Range("I102").Value = 4242.4634
$ws.Range("K102").Value = 4242.4634
$ws.Range("M102").Value = -2620.4634
# Row 132
$ws.Range("H132").Value = 4520.0293
$ws.Range("I132").Value = 4190.8
$ws.Range("J132").Value = 6989.25
$ws.Range("K132").Value = 12572.4
$ws.Range("L132").Value = 20967.75
$ws.Range("M132").Value = -10042.4
$ws.Range("N132").Value = -26027.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 6967.1113
$ws.Range("I122").Value = 5736.032
$ws.Range("K122").Value = 17208.096
$ws.Range("M122").Value = -14758.096
# Row 132
$ws.Range("H132").Value = 625397
$ws.Range("I132").Value = 1355868.2
$ws.Range("J132").Value = 7305.923
$ws.Range("K132").Value = 4067604.6
$ws.Range("L132").Value = 21917.769
$ws.Range("M132").Value = -4065074.6
$ws.Range("N132").Value = -26977.769
# Row 136
$ws.Range("H136").Value = 6782.5415
$ws.Range("I136").Value = 1446.2
$ws.Range("J136").Value = 15676.444
$ws.Range("K136").Value = 4338.6
$ws.Range("L136").Value = 47029.33199999999
$ws.Range("M136").Value = -1788.6
$ws.Range("N136").Value = -52129.33199999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value = 5937.75
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
# Row 36
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
# Row 45
$ws.Range("H45").Value = 11099.5
$ws.Range("J45").Value = 8319.4
$ws.Range("L45").Value = 8319.4
$ws.Range("N45").Value = -9301.4
# Row 81
$ws.Range("H81").Value = 12718.134
$ws.Range("I81").Value = 15798.363
$ws.Range("K81").Value = 31596.726
$ws.Range("M81").Value = -30535.726
# Row 84
$ws.Range("H84").Value = 12718.134
$ws.Range("I84").Value = 15798.363
$ws.Range("K84").Value = 157983.63
$ws.Range("M84").Value = -152679.63
# Row 96
$ws.Range("H96").Value = 3757.7
$ws.Range("I96").Value = 1654.1428
$ws.Range("J96").Value = 8666
$ws.Range("K96").Value = 1654.1428
$ws.Range("L96").Value = 8666
$ws.Range("M96").Value = -281.1428000000001
$ws.Range("N96").Value = -11412
# Row 132
$ws.Range("H132").Value = 7733.5
$ws.Range("I132").Value = 8096.2104
$ws.Range("J132").Value = 6517.353
$ws.Range("K132").Value = 24288.6312
$ws.Range("L132").Value = 19552.059
$ws.Range("M132").Value = -21758.6312
$ws.Range("N132").Value = -24612.059
# Row 136
$ws.Range("H136").Value = 483629
$ws.Range("I136").Value = 532832
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 1598496
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -1595946
$ws.Range("N136").Value = -29100
